# Quarterly financials update ("Doing Updates for Financials"):
# Two new quarterly columns (FY2018 Q4 ending 2018-12-31 and Q3 ending
# 2018-09-30) are inserted at the front of every data table on the sheet
# (Income Statement, Balance Sheet, Cash Flow). The existing quarter
# columns D:K shift right to F:M, and the two freshly inserted D:E
# columns are populated with the new quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank columns before column D -----------------------
# This shifts the existing data (columns D:K) to F:M and leaves two new
# blank columns at D:E, matching the formatting of the column that used
# to be D (now F) per Excel's native column-insert behaviour.
$ws.Range("D1:E1").EntireColumn.Insert()

# --- 2. Give the new D:E columns the same number formats as the data --
# columns immediately to their right (F:G), row by row, so dates keep
# the date format and amounts keep the accounting number format.
$ws.Range("F7:G102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 3. Also match the new columns' width/display to their neighbours -
$ws.Columns.Item("D").ColumnWidth = $ws.Columns.Item("F").ColumnWidth
$ws.Columns.Item("E").ColumnWidth = $ws.Columns.Item("G").ColumnWidth

# --- 4. Populate the new D:E columns with the new quarter's data ------
# Each entry is @(row, newDValue, newEValue); $null means the cell is
# left blank (spacer rows / section headers with no figures).
$rowData = @(
    @(7, 43465, 43373),
    @(8, 86500, 78000),
    @(9, 37000, 34500),
    @(10, 49500, 43500),
    @(11, $null, $null),
    @(12, 35500, 36500),
    @(13, 0, 0),
    @(14, 0, 0),
    @(15, 7800, 7600),
    @(16, $null, $null),
    @(17, 94800, 94100),
    @(18, -8300, -16100),
    @(19, $null, $null),
    @(20, -4800, 1400),
    @(21, 9700, 7800),
    @(22, 8300, 8200),
    @(23, -21400, -22900),
    @(24, 200, -200),
    @(25, 0, 0),
    @(26, -21600, -22700),
    @(27, -21600, -22700),
    @(28, 0, 0),
    @(29, 0, "NA"),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 4800, -1400),
    @(33, -21600, -22700),
    @(34, 0, 0),
    @(35, -21600, -22700),
    @(38, 43465, 43373),
    @(39, $null, $null),
    @(40, $null, $null),
    @(41, 172000, 172400),
    @(42, 235300, 222600),
    @(43, 61300, 48400),
    @(44, 33100, 34100),
    @(45, 9600, 8400),
    @(46, 511300, 485900),
    @(47, 0, 0),
    @(48, 70700, 67500),
    @(49, 284900, 301000),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 22900, 30900),
    @(53, 0, 0),
    @(54, 889900, 885400),
    @(55, $null, $null),
    @(56, $null, $null),
    @(57, 15900, 14300),
    @(58, "NA", 500),
    @(59, 48600, 37700),
    @(60, 64400, 52400),
    @(61, 447800, 441500),
    @(62, 10900, 14200),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 523200, 508100),
    @(67, $null, $null),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, -169900, -148300),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 366700, 377300),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, -21600, -22700),
    @(82, $null, $null),
    @(83, 22800, 22500),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 29800, 22300),
    @(90, $null, $null),
    @(91, -7800, -6200),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -25000, -5700),
    @(95, $null, $null),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, -5200, 700),
    @(101, 0, 0),
    @(102, -400, 17300)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    if ($null -ne $dVal) {
        $ws.Range("D$r").Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Range("E$r").Value = $eVal
    }
}

# --- 5. Keep the workbook's used range / calc mode consistent ---------
$ws.Range("A5:M102") | Out-Null
$excel.MultiThreadedCalculation = $false
